$d = $word.ActiveDocument

# Locate the end of the last paragraph of the document body: the one
# ending with "...It is always moving from parent to child component."
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "It is always moving from parent to child component.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor sentence not found; cannot place the new Q&A paragraphs."
}

$tail = $anchor.Duplicate
$tail.Collapse(0)

# --- New paragraph 1: the bold, numbered question ------------------------
$tail.InsertParagraphAfter() | Out-Null
$qPara = $d.Paragraphs.Last
$qRange = $qPara.Range

$qXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="202122"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="202122"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>What is the use of a super keyword in React?</w:t>
  </w:r>
</w:p>
"@

$qRange.InsertXML($qXml) | Out-Null

# --- New paragraph 2: the (non-bold) answer text --------------------------
$qPara2 = $d.Paragraphs.Last
$qRange2 = $qPara2.Range.Duplicate
$qRange2.Collapse(0)
$qRange2.InsertParagraphAfter() | Out-Null

$aPara = $d.Paragraphs.Last
$aRange = $aPara.Range

$aXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:ind w:left="360"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:color w:val="202122"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:color w:val="202122"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>The super keyword helps you to access and call functions on an object&#8217;s parent.</w:t>
  </w:r>
</w:p>
"@

$aRange.InsertXML($aXml) | Out-Null

Write-Output "Inserted the super-keyword Q&A paragraphs."
